$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date serial values from 45184 to 45186
# for all data rows (rows 2 through 23).
$ws.Range("C2:C23").Value = 45186
